$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows covering 2021-05-28 through 2021-06-28 (update "fino a 28/06 incluso").
# Columns: A=date serial, B=nuovi pos., C=somma mobile 7gg., D=somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44344, 0, 6, 38.75217981011431),
    @(44345, 0, 2, 12.91739327003811),
    @(44346, 0, 2, 12.91739327003811),
    @(44347, 0, 0, 0),
    @(44348, 0, 0, 0),
    @(44349, 0, 0, 0),
    @(44350, 0, 0, 0),
    @(44351, 1, 1, 6.458696635019054),
    @(44352, 0, 1, 6.458696635019054),
    @(44353, 0, 1, 6.458696635019054),
    @(44354, 0, 1, 6.458696635019054),
    @(44355, 1, 2, 12.91739327003811),
    @(44356, 1, 3, 19.37608990505716),
    @(44357, 0, 3, 19.37608990505716),
    @(44358, 1, 3, 19.37608990505716),
    @(44359, 1, 4, 25.83478654007622),
    @(44360, 0, 4, 25.83478654007622),
    @(44361, 0, 4, 25.83478654007622),
    @(44362, 0, 3, 19.37608990505716),
    @(44363, 0, 2, 12.91739327003811),
    @(44364, 0, 2, 12.91739327003811),
    @(44365, 0, 1, 6.458696635019054),
    @(44366, 1, 1, 6.458696635019054),
    @(44367, 0, 1, 6.458696635019054),
    @(44368, 0, 1, 6.458696635019054),
    @(44369, 0, 1, 6.458696635019054),
    @(44370, 0, 1, 6.458696635019054),
    @(44371, 0, 1, 6.458696635019054),
    @(44372, 0, 1, 6.458696635019054),
    @(44373, 0, 0, 0),
    @(44374, 0, 0, 0),
    @(44375, 0, 0, 0)
)

$startRow = 270
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
$endRow = $r - 1

# Match the date-column formatting/border/style already used by the existing data (style index "2").
$ws.Range("A269").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

Write-Host "Added rows $startRow to $endRow"
